# Add the new book "In Defense of Elitism" as the next row in the
# "Completed" reading-list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")
$ws.Activate()

# New row goes right after the current last row (row 40).
$newRow = 41

# Carry the date formatting (style) down from the row above first, so the
# new date cells use the same existing date style instead of Excel
# creating a brand new number format.
$ws.Range("C40:D40").Copy()
$ws.Range("C41:D41").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values in the order they were entered.
$ws.Cells.Item($newRow, 7).Value = "7 Hours 19 Mins"
$ws.Cells.Item($newRow, 1).Value = "In Defense of Elitism"
$ws.Cells.Item($newRow, 2).Value = "Joel Stein"
$ws.Cells.Item($newRow, 3).Value = "3/20/2020"
$ws.Cells.Item($newRow, 4).Value = "3/22/2020"
$ws.Cells.Item($newRow, 5).Value = "elitism;academics;populism;politics"
$ws.Cells.Item($newRow, 6).Value = "Audio"

# Leave the selection where Excel would land after typing across the row.
$ws.Range("C42").Select()
